$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("A1").Value = "TestValue"
